# overview_testcases.xlsx - "removed xlsx cached stuff"
#
# Core data edit: I39 on the "overview_testcases" sheet changes from 128 to
# 512. I40:I48 are `=$I$39` formulas (and K39:K48 are the
# MAX/AVERAGE-derived ratio in column K), so writing the new literal lets
# the engine's auto-recalc cascade the dependent formula cells/cached
# values for us.
#
# The diff also shows the workbook re-opening with the first sheet
# ("overview_testcases") focused/active (instead of "overhead"), with a
# new selection/scroll position on that sheet, and "overhead" losing its
# tabSelected flag.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("overview_testcases")
$wsOverhead = $wb.Worksheets.Item("overhead")

# --- the actual data change -------------------------------------------------
$wsOverview.Range("I39").Value = 512

# --- view / selection state --------------------------------------------------
# Make "overview_testcases" the active sheet/tab (was "overhead").
$wsOverview.Activate() | Out-Null

$win = $excel.ActiveWindow
$win.Zoom = 80

# Scroll so row 6 is at the top, and move the selection to M43 (matches the
# new <selection activeCell="M43" sqref="M43"/> on the sheet).
$win.ScrollRow = 6
$win.ScrollColumn = 1
$wsOverview.Range("M43").Select() | Out-Null

# "overhead" sheet keeps its own prior selection (G15); just make sure it is
# no longer the tab-selected sheet, which Activate() above already took care
# of.
$wsOverhead.Range("G15").Select() | Out-Null

# Re-activate overview_testcases so it remains the visible/active sheet and
# M43 remains the active selection on it.
$wsOverview.Activate() | Out-Null
